$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Create")

$ws1.Range("C3").Value = "'4356"
$ws1.Range("E3").Value = "'4728"
$ws1.Range("F3").Value = "'Rest"

$ws1.Range("E4").Value = "'5267"
$ws1.Range("F4").Value = "'User"

$ws1.Range("F5").Value = "'Data"

$ws4 = $wb.Worksheets.Item("Queries")
$ws4.Range("A2").Value = "Select IMU_MENU_ID as 'Menu ID'
  ,IVN_VDN_NUM as VDN
  ,IVN_VDN_DESC as 'VDN Description'
  ,IVN_VIP_VDN_NUM as 'VIP VDN'
  ,IVN_MENU_OPT as 'Option'
  ,IVN_SESS_DNIS as DNIS
  FROM [IVR_VDN_NUM] Order By IVN_VDN_DESC Asc;"
